$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("BD3").Value = 151
$ws.Range("O5").Value = 1.5
$ws.Range("P5").Value = 2.5
$ws.Range("O10").Value = 1.13
$ws.Range("U10").Value = 2.28
$ws.Range("V10").Value = 1.6
$ws.Range("O16").Value = 1.2
$ws.Range("P16").Value = 4.33
$ws.Range("Q16").Value = 1.65
$ws.Range("R16").Value = 2.2
$ws.Range("Q17").Value = 1.88
$ws.Range("R17").Value = 1.98
$ws.Range("J22").Value = 1.8
$ws.Range("K22").Value = 2.88
$ws.Range("Q22").Value = 1.33
$ws.Range("R22").Value = 3.4
$ws.Range("M27").Value = 1.08
$ws.Range("N27").Value = 8
$ws.Range("M30").Value = 1.06
$ws.Range("N30").Value = 10
$ws.Range("O30").Value = 1.29
$ws.Range("Q30").Value = 2
$ws.Range("R30").Value = 1.85
$ws.Range("M32").Value = 1.03
$ws.Range("O32").Value = 1.14
